$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Delete rows from the bottom up so earlier row numbers stay valid references.
# Original rows 32 & 33 (45361.99999999999/10 and 45368.99999999999/20) are removed entirely.
$ws1.Range("A32:A33").EntireRow.Delete()

# Original row 13 (45137.99999999999/290) is removed entirely, shifting later rows up.
$ws1.Range("A13").EntireRow.Delete()

# --- Sheet "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Original row 12 (45382.99999999999/30) is removed entirely.
$ws2.Range("A12").EntireRow.Delete()

# Row 5's requested quantity changes from 400 to 110.
$ws2.Range("B5").Value = 110
